$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 34: 2940803 / LE CHAT ROSE 1L / 10 / 377 ---

# A34: Item ID stored as TEXT ("2940803") rather than a number, matching
# the source data. We build it via a helper formula cell that returns text,
# copy it, and paste-special the VALUE into A34 so the literal becomes a
# plain shared string (no residual formula, no quote-prefix artifacts).
$helper = $ws.Range("Z1")
$helper.Formula = '="2940803"'
$helper.Copy()
$ws.Range("A34").PasteSpecial(-4163)
$helper.ClearContents()

$ws.Range("A34").Borders.LineStyle = 1
$ws.Range("A34").HorizontalAlignment = -4131

$ws.Range("B34").Borders.LineStyle = 1
$ws.Range("B34").Value = "LE CHAT ROSE 1L"

$ws.Range("C34").Borders.LineStyle = 1
$ws.Range("C34").Value = 10

$ws.Range("D34").Borders.LineStyle = 1
$ws.Range("D34").Value = 377

# --- Row 35: 2940804 / LE CHAT power gel 4L / 4 / 1190 ---

$ws.Range("A35").Borders.LineStyle = 1
$ws.Range("A35").HorizontalAlignment = -4131
$ws.Range("A35").Value = 2940804

$ws.Range("B35").Borders.LineStyle = 1
$ws.Range("B35").Value = "LE CHAT power gel 4L"

$ws.Range("C35").Borders.LineStyle = 1
$ws.Range("C35").Value = 4

$ws.Range("D35").Borders.LineStyle = 1
$ws.Range("D35").Value = 1190

# --- View state: scroll/select like the saved workbook ---
$ws.Range("A34:D35").Select()
$excel.ActiveWindow.ScrollRow = 25
